$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of consultation data appended below the existing header/record row.
$ws.Range("A2").Value = "js"
$ws.Range("B2").Value = "+91 72489 46823"
$ws.Range("C2").Value = "Headache"
$ws.Range("D2").Value = "Dr. Joseph King: Plastic Surgeon"

# Force text format so the date/time strings are preserved literally,
# matching the source data instead of being auto-converted to date/time serials.
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "2024-09-26"
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "15:24"
